$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lecture_Schedule_SS2022")

# Shift the lecture-date column (E) forward by one year (2021 -> 2022).
# E15 is the "anchor" literal date for the odd-row chain; every formula
# cell below it (E17, E19, ... E43, each "=<prior E cell>+7") recalculates
# automatically once E15 changes.
$ws.Range("E15").Value2 = 44672

# E44 used to hold the shared formula "=E42+7"; it is now overwritten with
# a manually-entered literal date value instead.
$ws.Range("E44").Value2 = 44774

# Update the view: scroll back up (drops the saved topLeftCell="A12") and
# move the selection from H36 to E16.
$ws.Activate() | Out-Null
$ws.Range("E16").Select() | Out-Null
